$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1154.1072
$ws.Range("I121").Value = 1500
$ws.Range("J121").Value = 1141.2963
$ws.Range("K121").Value = 4500
$ws.Range("L121").Value = 3423.8889
$ws.Range("M121").Value = -2753
$ws.Range("N121").Value = -6917.8889

$ws.Range("H132").Value = 30836532
$ws.Range("I132").Value = 42695580
$ws.Range("J132").Value = 3010.5
$ws.Range("K132").Value = 128086740
$ws.Range("L132").Value = 9031.5
$ws.Range("M132").Value = -128084210
$ws.Range("N132").Value = -14091.5

$ws.Range("H137").Value = 166459.94
$ws.Range("I137").Value = 259750.95
$ws.Range("J137").Value = 1406.5769
$ws.Range("K137").Value = 779252.8500000001
$ws.Range("L137").Value = 4219.7307
$ws.Range("M137").Value = -776702.8500000001
$ws.Range("N137").Value = -9319.7307

$ws.Range("H138").Value = 3608.5
$ws.Range("I138").Value = 5876.7188
$ws.Range("J138").Value = 2096.3542
$ws.Range("K138").Value = 17630.1564
$ws.Range("L138").Value = 6289.062600000001
$ws.Range("M138").Value = -12490.1564
$ws.Range("N138").Value = -16569.0626

$ws.Range("H141").Value = 11956.363
$ws.Range("I141").Value = 16753.334
$ws.Range("J141").Value = 6200
$ws.Range("K141").Value = 50260.00199999999
$ws.Range("L141").Value = 18600
$ws.Range("M141").Value = -45080.00199999999
$ws.Range("N141").Value = -28960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 7574.2856
$ws.Range("I37").Value = 5034
$ws.Range("J37").Value = 7997.6665
$ws.Range("K37").Value = 5034
$ws.Range("L37").Value = 7997.6665
$ws.Range("M37").Value = -4761
$ws.Range("N37").Value = -8543.666499999999

$ws.Range("H55").Value = 28744.75
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 34993
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 34993
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -35623

$ws.Range("H61").Value = 4133.3335
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 4133.3335
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4133.3335
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -4557.3335

$ws.Range("H74").Value = 2162.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2162.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2162.8
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -3910.8

$ws.Range("H77").Value = 2162.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2162.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 10814
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -19550

$ws.Range("H80").Value = 30274
$ws.Range("J80").Value = 30274
$ws.Range("L80").Value = 30274
$ws.Range("N80").Value = -32270

$ws.Range("H83").Value = 30274
$ws.Range("J83").Value = 30274
$ws.Range("L83").Value = 90822
$ws.Range("N83").Value = -100806

$ws.Range("H132").Value = 2605887.8
$ws.Range("I132").Value = 4311355.5
$ws.Range("J132").Value = 2805.7896
$ws.Range("K132").Value = 12934066.5
$ws.Range("L132").Value = 8417.3688
$ws.Range("M132").Value = -12931536.5
$ws.Range("N132").Value = -13477.3688

$ws.Range("H136").Value = 4133.3335
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4133.3335
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 12400.0005
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -17500.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31249.75
$ws.Range("I20").Value = 36666.332
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 36666.332
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -36419.332
$ws.Range("N20").Value = -15494

$ws.Range("H134").Value = 3109.4546
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 3109.4546
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 9328.363799999999
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -14398.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12992062
$ws.Range("I31").Value = 15152401
$ws.Range("J31").Value = 30027.908
$ws.Range("K31").Value = 15152401
$ws.Range("L31").Value = 30027.908
$ws.Range("M31").Value = -15152106
$ws.Range("N31").Value = -30617.908

$ws.Range("H34").Value = 12992062
$ws.Range("I34").Value = 15152401
$ws.Range("J34").Value = 30027.908
$ws.Range("K34").Value = 15152401
$ws.Range("L34").Value = 30027.908
$ws.Range("M34").Value = -15152199
$ws.Range("N34").Value = -30431.908

$ws.Range("H58").Value = 2286584.8
$ws.Range("I58").Value = 4232504.5
$ws.Range("J58").Value = 5161.1035
$ws.Range("K58").Value = 4232504.5
$ws.Range("L58").Value = 5161.1035
$ws.Range("M58").Value = -4232301.5
$ws.Range("N58").Value = -5567.1035

$ws.Range("H94").Value = 2811.2
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2811.2
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2811.2
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -3713.2

$ws.Range("H132").Value = 8919.883
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 8919.883
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 26759.649
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -31819.649

$ws.Range("H134").Value = 50482830
$ws.Range("I134").Value = 104168780
$ws.Range("J134").Value = 4466306.5
$ws.Range("K134").Value = 312506340
$ws.Range("L134").Value = 13398919.5
$ws.Range("M134").Value = -312503805
$ws.Range("N134").Value = -13403989.5

$ws.Range("H136").Value = 2286584.8
$ws.Range("I136").Value = 4232504.5
$ws.Range("J136").Value = 5161.1035
$ws.Range("K136").Value = 12697513.5
$ws.Range("L136").Value = 15483.3105
$ws.Range("M136").Value = -12694963.5
$ws.Range("N136").Value = -20583.3105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19253062
$ws.Range("I132").Value = 31282644
$ws.Range("J132").Value = 5730.4
$ws.Range("K132").Value = 93847932
$ws.Range("L132").Value = 17191.2
$ws.Range("M132").Value = -93845402
$ws.Range("N132").Value = -22251.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2818770.5
$ws.Range("I132").Value = 7144900
$ws.Range("J132").Value = 1756.0233
$ws.Range("K132").Value = 21434700
$ws.Range("L132").Value = 5268.0699
$ws.Range("M132").Value = -21432170
$ws.Range("N132").Value = -10328.0699

$ws.Range("H136").Value = 3349.3872
$ws.Range("I136").Value = 4086.2307
$ws.Range("J136").Value = 2099.9565
$ws.Range("K136").Value = 12258.6921
$ws.Range("L136").Value = 6299.869499999999
$ws.Range("M136").Value = -9708.6921
$ws.Range("N136").Value = -11399.8695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13490935
$ws.Range("I132").Value = 6484627.5
$ws.Range("K132").Value = 19453882.5
$ws.Range("M132").Value = -19451352.5

$ws.Range("H136").Value = 16483333
$ws.Range("I136").Value = 8305014.5
$ws.Range("J136").Value = 55557520
$ws.Range("K136").Value = 24915043.5
$ws.Range("L136").Value = 166672560
$ws.Range("M136").Value = -24912493.5
$ws.Range("N136").Value = -166677660
